$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 and 4 (bottom-up to keep indices stable)
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# Delete columns E through I (bottom-up to keep indices stable)
$ws.Columns.Item(9).Delete()
$ws.Columns.Item(8).Delete()
$ws.Columns.Item(7).Delete()
$ws.Columns.Item(6).Delete()
$ws.Columns.Item(5).Delete()

# Update header row
$ws.Range("A1").Value = "Qtd_Nós"
$ws.Range("B1").Value = "Ativos"
$ws.Range("D1").Value = "Tempo"

# Update data row
$ws.Range("A2").Value = 42
$ws.Range("B2").Value = 25
$ws.Range("D2").Value = 1880.658307790756
